$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new product rows (139-175) for the "letreros-banners" category ---
# (products_complete.xlsx template-rendering fix: new category data appended)

$ws.Cells.Item(139,1).Value = "banderas-gota-de-agua"
$ws.Cells.Item(139,2).Value = "Banderas Gota De Agua"
$ws.Cells.Item(139,3).Value = "letreros-banners"
$ws.Cells.Item(139,4).Value = "banderas"
$ws.Cells.Item(139,5).Value = "BAN-001"
$ws.Cells.Item(139,7).Value = "/media/product_images/letreros_banners/banderas/banderas-gota-de-agua.jpg"

$ws.Cells.Item(140,1).Value = "banderas-marinero"
$ws.Cells.Item(140,2).Value = "Banderas Marinero"
$ws.Cells.Item(140,3).Value = "letreros-banners"
$ws.Cells.Item(140,4).Value = "banderas"
$ws.Cells.Item(140,5).Value = "BAN-002"
$ws.Cells.Item(140,7).Value = "/media/product_images/letreros_banners/banderas/banderas-marinero.jpg"

$ws.Cells.Item(141,1).Value = "banderas-pared"
$ws.Cells.Item(141,2).Value = "Banderas Pared"
$ws.Cells.Item(141,3).Value = "letreros-banners"
$ws.Cells.Item(141,4).Value = "banderas"
$ws.Cells.Item(141,5).Value = "BAN-003"
$ws.Cells.Item(141,7).Value = "/media/product_images/letreros_banners/banderas/banderas-pared.jpg"

$ws.Cells.Item(142,1).Value = "banderas-pluma"
$ws.Cells.Item(142,2).Value = "Banderas Pluma"
$ws.Cells.Item(142,3).Value = "letreros-banners"
$ws.Cells.Item(142,4).Value = "banderas"
$ws.Cells.Item(142,5).Value = "BAN-004"
$ws.Cells.Item(142,7).Value = "/media/product_images/letreros_banners/banderas/banderas-pluma.jpg"

$ws.Cells.Item(143,1).Value = "banderas-polyester"
$ws.Cells.Item(143,2).Value = "Banderas Polyester"
$ws.Cells.Item(143,3).Value = "letreros-banners"
$ws.Cells.Item(143,4).Value = "banderas"
$ws.Cells.Item(143,5).Value = "BAN-005"
$ws.Cells.Item(143,7).Value = "/media/product_images/letreros_banners/banderas/banderas-polyester.jpg"

$ws.Cells.Item(144,1).Value = "banderas-rectangular"
$ws.Cells.Item(144,2).Value = "Banderas Rectangular"
$ws.Cells.Item(144,3).Value = "letreros-banners"
$ws.Cells.Item(144,4).Value = "banderas"
$ws.Cells.Item(144,5).Value = "BAN-006"
$ws.Cells.Item(144,7).Value = "/media/product_images/letreros_banners/banderas/banderas-rectangular.jpg"

$ws.Cells.Item(145,1).Value = "banderines"
$ws.Cells.Item(145,2).Value = "Banderines"
$ws.Cells.Item(145,3).Value = "letreros-banners"
$ws.Cells.Item(145,4).Value = "banderas"
$ws.Cells.Item(145,5).Value = "BAN-007"
$ws.Cells.Item(145,7).Value = "/media/product_images/letreros_banners/banderas/banderines.jpg"

$ws.Cells.Item(146,1).Value = "banners-malla"
$ws.Cells.Item(146,2).Value = "Banners Malla"
$ws.Cells.Item(146,3).Value = "letreros-banners"
$ws.Cells.Item(146,4).Value = "banners"
$ws.Cells.Item(146,5).Value = "BAN-008"
$ws.Cells.Item(146,7).Value = "/media/product_images/letreros_banners/banners/banners-malla.jpg"

$ws.Cells.Item(147,1).Value = "banners-polyester"
$ws.Cells.Item(147,2).Value = "Banners Polyester"
$ws.Cells.Item(147,3).Value = "letreros-banners"
$ws.Cells.Item(147,4).Value = "banners"
$ws.Cells.Item(147,5).Value = "BAN-009"
$ws.Cells.Item(147,7).Value = "/media/product_images/letreros_banners/banners/banners-polyester.jpg"

$ws.Cells.Item(148,1).Value = "banners-postes"
$ws.Cells.Item(148,2).Value = "Banners Postes"
$ws.Cells.Item(148,3).Value = "letreros-banners"
$ws.Cells.Item(148,4).Value = "banners"
$ws.Cells.Item(148,5).Value = "BAN-010"
$ws.Cells.Item(148,7).Value = "/media/product_images/letreros_banners/banners/banners-postes.jpg"

$ws.Cells.Item(149,1).Value = "banners-repaso-repeticion"
$ws.Cells.Item(149,2).Value = "Banners Repaso Repeticion"
$ws.Cells.Item(149,3).Value = "letreros-banners"
$ws.Cells.Item(149,4).Value = "banners"
$ws.Cells.Item(149,5).Value = "BAN-011"
$ws.Cells.Item(149,7).Value = "/media/product_images/letreros_banners/banners/banners-repaso-repeticion.jpg"

$ws.Cells.Item(150,1).Value = "banners-retractables"
$ws.Cells.Item(150,2).Value = "Banners Retractables"
$ws.Cells.Item(150,3).Value = "letreros-banners"
$ws.Cells.Item(150,4).Value = "banners"
$ws.Cells.Item(150,5).Value = "BAN-012"
$ws.Cells.Item(150,7).Value = "/media/product_images/letreros_banners/banners/banners-retractables.jpg"

$ws.Cells.Item(151,1).Value = "banners-tensados"
$ws.Cells.Item(151,2).Value = "Banners Tensados"
$ws.Cells.Item(151,3).Value = "letreros-banners"
$ws.Cells.Item(151,4).Value = "banners"
$ws.Cells.Item(151,5).Value = "BAN-013"
$ws.Cells.Item(151,7).Value = "/media/product_images/letreros_banners/banners/banners-tensados.jpg"

$ws.Cells.Item(152,1).Value = "banners-vinyl"
$ws.Cells.Item(152,2).Value = "Banners Vinyl"
$ws.Cells.Item(152,3).Value = "letreros-banners"
$ws.Cells.Item(152,4).Value = "banners"
$ws.Cells.Item(152,5).Value = "BAN-014"
$ws.Cells.Item(152,7).Value = "/media/product_images/letreros_banners/banners/banners-vinyl.jpg"

$ws.Cells.Item(153,1).Value = "banners-x"
$ws.Cells.Item(153,2).Value = "Banners X"
$ws.Cells.Item(153,3).Value = "letreros-banners"
$ws.Cells.Item(153,4).Value = "banners"
$ws.Cells.Item(153,5).Value = "BAN-015"
$ws.Cells.Item(153,7).Value = "/media/product_images/letreros_banners/banners/banners-x.jpg"

$ws.Cells.Item(154,1).Value = "letreros-0001"
$ws.Cells.Item(154,2).Value = "Letreros 0001"
$ws.Cells.Item(154,3).Value = "letreros-banners"
$ws.Cells.Item(154,4).Value = "letreros"
$ws.Cells.Item(154,5).Value = "LE-001"
$ws.Cells.Item(154,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0001.jpg"

$ws.Cells.Item(155,1).Value = "letreros-0002"
$ws.Cells.Item(155,2).Value = "Letreros 0002"
$ws.Cells.Item(155,3).Value = "letreros-banners"
$ws.Cells.Item(155,4).Value = "letreros"
$ws.Cells.Item(155,5).Value = "LE-002"
$ws.Cells.Item(155,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0002.jpg"

$ws.Cells.Item(156,1).Value = "letreros-0003"
$ws.Cells.Item(156,2).Value = "Letreros 0003"
$ws.Cells.Item(156,3).Value = "letreros-banners"
$ws.Cells.Item(156,4).Value = "letreros"
$ws.Cells.Item(156,5).Value = "LE-003"
$ws.Cells.Item(156,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0003.jpg"

$ws.Cells.Item(157,1).Value = "letreros-0004"
$ws.Cells.Item(157,2).Value = "Letreros 0004"
$ws.Cells.Item(157,3).Value = "letreros-banners"
$ws.Cells.Item(157,4).Value = "letreros"
$ws.Cells.Item(157,5).Value = "LE-004"
$ws.Cells.Item(157,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0004.jpg"

$ws.Cells.Item(158,1).Value = "letreros-0005"
$ws.Cells.Item(158,2).Value = "Letreros 0005"
$ws.Cells.Item(158,3).Value = "letreros-banners"
$ws.Cells.Item(158,4).Value = "letreros"
$ws.Cells.Item(158,5).Value = "LE-005"
$ws.Cells.Item(158,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0005.jpg"

$ws.Cells.Item(159,1).Value = "letreros-0006"
$ws.Cells.Item(159,2).Value = "Letreros 0006"
$ws.Cells.Item(159,3).Value = "letreros-banners"
$ws.Cells.Item(159,4).Value = "letreros"
$ws.Cells.Item(159,5).Value = "LE-006"
$ws.Cells.Item(159,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0006.jpg"

$ws.Cells.Item(160,1).Value = "letreros-0007"
$ws.Cells.Item(160,2).Value = "Letreros 0007"
$ws.Cells.Item(160,3).Value = "letreros-banners"
$ws.Cells.Item(160,4).Value = "letreros"
$ws.Cells.Item(160,5).Value = "LE-007"
$ws.Cells.Item(160,7).Value = "/media/product_images/letreros_banners/letreros/letreros-0007.jpg"

$ws.Cells.Item(161,1).Value = "posters-0001"
$ws.Cells.Item(161,2).Value = "Posters 0001"
$ws.Cells.Item(161,3).Value = "letreros-banners"
$ws.Cells.Item(161,4).Value = "posters"
$ws.Cells.Item(161,5).Value = "LE-008"
$ws.Cells.Item(161,7).Value = "/media/product_images/letreros_banners/posters/posters-0001.jpg"

$ws.Cells.Item(162,1).Value = "posters-0002"
$ws.Cells.Item(162,2).Value = "Posters 0002"
$ws.Cells.Item(162,3).Value = "letreros-banners"
$ws.Cells.Item(162,4).Value = "posters"
$ws.Cells.Item(162,5).Value = "LE-009"
$ws.Cells.Item(162,7).Value = "/media/product_images/letreros_banners/posters/posters-0002.jpg"

$ws.Cells.Item(163,1).Value = "posters-0004"
$ws.Cells.Item(163,2).Value = "Posters 0004"
$ws.Cells.Item(163,3).Value = "letreros-banners"
$ws.Cells.Item(163,4).Value = "posters"
$ws.Cells.Item(163,5).Value = "LE-010"
$ws.Cells.Item(163,7).Value = "/media/product_images/letreros_banners/posters/posters-0004.jpg"

$ws.Cells.Item(164,1).Value = "posters-0005"
$ws.Cells.Item(164,2).Value = "Posters 0005"
$ws.Cells.Item(164,3).Value = "letreros-banners"
$ws.Cells.Item(164,4).Value = "posters"
$ws.Cells.Item(164,5).Value = "LE-011"
$ws.Cells.Item(164,7).Value = "/media/product_images/letreros_banners/posters/posters-0005.jpg"

$ws.Cells.Item(165,1).Value = "posters-0006"
$ws.Cells.Item(165,2).Value = "Posters 0006"
$ws.Cells.Item(165,3).Value = "letreros-banners"
$ws.Cells.Item(165,4).Value = "posters"
$ws.Cells.Item(165,5).Value = "LE-012"
$ws.Cells.Item(165,7).Value = "/media/product_images/letreros_banners/posters/posters-0006.jpg"

$ws.Cells.Item(166,1).Value = "posters-0007"
$ws.Cells.Item(166,2).Value = "Posters 0007"
$ws.Cells.Item(166,3).Value = "letreros-banners"
$ws.Cells.Item(166,4).Value = "posters"
$ws.Cells.Item(166,5).Value = "LE-013"
$ws.Cells.Item(166,7).Value = "/media/product_images/letreros_banners/posters/posters-0007.jpg"

$ws.Cells.Item(167,1).Value = "posters-0008"
$ws.Cells.Item(167,2).Value = "Posters 0008"
$ws.Cells.Item(167,3).Value = "letreros-banners"
$ws.Cells.Item(167,4).Value = "posters"
$ws.Cells.Item(167,5).Value = "LE-014"
$ws.Cells.Item(167,7).Value = "/media/product_images/letreros_banners/posters/posters-0008.jpg"

$ws.Cells.Item(168,1).Value = "posters-0009"
$ws.Cells.Item(168,2).Value = "Posters 0009"
$ws.Cells.Item(168,3).Value = "letreros-banners"
$ws.Cells.Item(168,4).Value = "posters"
$ws.Cells.Item(168,5).Value = "LE-015"
$ws.Cells.Item(168,7).Value = "/media/product_images/letreros_banners/posters/posters-0009.jpg"

$ws.Cells.Item(169,1).Value = "posters-exteriores"
$ws.Cells.Item(169,2).Value = "Posters Exteriores"
$ws.Cells.Item(169,3).Value = "letreros-banners"
$ws.Cells.Item(169,4).Value = "posters"
$ws.Cells.Item(169,5).Value = "LE-016"
$ws.Cells.Item(169,7).Value = "/media/product_images/letreros_banners/posters/posters-exteriores.jpg"

$ws.Cells.Item(170,1).Value = "publicidad-autos-0001"
$ws.Cells.Item(170,2).Value = "Publicidad Autos 0001"
$ws.Cells.Item(170,3).Value = "letreros-banners"
$ws.Cells.Item(170,4).Value = "publicidad-autos"
$ws.Cells.Item(170,5).Value = "PA-001"
$ws.Cells.Item(170,7).Value = "/media/product_images/letreros_banners/publicidad_autos/publicidad-autos-0001.jpg"

$ws.Cells.Item(171,1).Value = "publicidad-autos-0002"
$ws.Cells.Item(171,2).Value = "Publicidad Autos 0002"
$ws.Cells.Item(171,3).Value = "letreros-banners"
$ws.Cells.Item(171,4).Value = "publicidad-autos"
$ws.Cells.Item(171,5).Value = "PA-002"
$ws.Cells.Item(171,7).Value = "/media/product_images/letreros_banners/publicidad_autos/publicidad-autos-0002.jpg"

$ws.Cells.Item(172,1).Value = "publicidad-autos-0003"
$ws.Cells.Item(172,2).Value = "Publicidad Autos 0003"
$ws.Cells.Item(172,3).Value = "letreros-banners"
$ws.Cells.Item(172,4).Value = "publicidad-autos"
$ws.Cells.Item(172,5).Value = "PA-003"
$ws.Cells.Item(172,7).Value = "/media/product_images/letreros_banners/publicidad_autos/publicidad-autos-0003.jpg"

$ws.Cells.Item(173,1).Value = "publicidad-autos-0004"
$ws.Cells.Item(173,2).Value = "Publicidad Autos 0004"
$ws.Cells.Item(173,3).Value = "letreros-banners"
$ws.Cells.Item(173,4).Value = "publicidad-autos"
$ws.Cells.Item(173,5).Value = "PA-004"
$ws.Cells.Item(173,7).Value = "/media/product_images/letreros_banners/publicidad_autos/publicidad-autos-0004.jpg"

$ws.Cells.Item(174,1).Value = "publicidad-autos-0005"
$ws.Cells.Item(174,2).Value = "Publicidad Autos 0005"
$ws.Cells.Item(174,3).Value = "letreros-banners"
$ws.Cells.Item(174,4).Value = "publicidad-autos"
$ws.Cells.Item(174,5).Value = "PA-005"
$ws.Cells.Item(174,7).Value = "/media/product_images/letreros_banners/publicidad_autos/publicidad-autos-0005.jpg"

$ws.Cells.Item(175,1).Value = "publicidad-autos-0006"
$ws.Cells.Item(175,2).Value = "Publicidad Autos 0006"
$ws.Cells.Item(175,3).Value = "letreros-banners"
$ws.Cells.Item(175,4).Value = "publicidad-autos"
$ws.Cells.Item(175,5).Value = "PA-006"
$ws.Cells.Item(175,7).Value = "/media/product_images/letreros_banners/publicidad_autos/publicidad-autos-0006.jpg"

# --- Adjust column widths (A-E) to fit the newly-added columns/content ---
$ws.Columns.Item(1).ColumnWidth = 37.6666666666667
$ws.Columns.Item(2).ColumnWidth = 34.1666666666667
$ws.Columns.Item(3).ColumnWidth = 18.3333333333333
$ws.Columns.Item(4).ColumnWidth = 19
$ws.Columns.Item(5).ColumnWidth = 31.6666666666667

# --- Apply the AutoFilter on the full data range, filtering to the new category ---
# (this also hides every pre-existing row 2-138 whose category_slug does not match,
#  and keeps the new 139-175 rows visible, matching the authored workbook state)
$ws.Range("A1:O175").AutoFilter(3, @("letreros-banners"), 7)

# --- Keep the _FilterDatabase defined name in sync with the new filtered range ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$O`$175"
    }
}

# --- Reset the view: selection on B1, scrolled back to the top-left ---
$ws.Range("B1").Select()
